# Fill in the training records for "Mohamed Khairy Elsayed Elmasry"
# (DSS1429, DSS1430, DSS1432, DSS1433, DSS1435, DSS1436, DSS1437, DSS1438)
# into the previously-blank placeholder rows 430, 431, 433, 434, 436-439,
# and fix up the style of A435 (was using the "date" flavoured style,
# now uses the plain text one, matching its siblings A432/B435/C435).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Stamp each target row with the exact existing cell styles that the
#    workbook already uses for this kind of record (name/ID rows), by
#    copying format from donor cells that already carry those styles.
#    This reuses the existing cellXfs entries instead of synthesizing
#    new ones.
# ---------------------------------------------------------------------
$idRows = @(430, 431, 433, 434, 436, 437, 438, 439)

foreach ($r in $idRows) {
    $ws.Range("C99").Copy()
    $ws.Range("A$r").PasteSpecial(-4122)
    $ws.Range("B$r").PasteSpecial(-4122)

    $ws.Range("C116").Copy()
    $ws.Range("C$r").PasteSpecial(-4122)

    $ws.Range("D211").Copy()
    $ws.Range("D$r").PasteSpecial(-4122)

    $ws.Range("E2").Copy()
    $ws.Range("E$r").PasteSpecial(-4122)
}

# A435 switches from the "mmm-yy date" flavoured style to the plain one
# (matching A432's style), value/content unchanged.
$ws.Range("C102").Copy()
$ws.Range("A435").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2) Write the values. The employee name is entered into column B first
#    (row 430) so it lands in the shared-string table before the first
#    new DSS id, matching how the data was actually typed in.
# ---------------------------------------------------------------------
$ws.Range("B430").Value = "Mohamed Khairy Elsayed Elmasry"

$ws.Range("A430").Value = "DSS1429"
$ws.Range("C430").Value = "30 Hours Construction Safety & Health"
$ws.Range("D430").Value = 45566
$ws.Range("E430").Value = 1

$ws.Range("A431").Value = "DSS1430"
$ws.Range("B431").Value = "Mohamed Khairy Elsayed Elmasry"
$ws.Range("C431").Value = "30 Hours G. Industry Safety & Health"
$ws.Range("D431").Value = 45566
$ws.Range("E431").Value = 1

$ws.Range("A433").Value = "DSS1432"
$ws.Range("B433").Value = "Mohamed Khairy Elsayed Elmasry"
$ws.Range("C433").Value = "Electrical Safety "
$ws.Range("D433").Value = 45566
$ws.Range("E433").Value = 1

$ws.Range("A434").Value = "DSS1433"
$ws.Range("B434").Value = "Mohamed Khairy Elsayed Elmasry"
$ws.Range("C434").Value = "Fire Marshal"
$ws.Range("D434").Value = 45566
$ws.Range("E434").Value = 1

$ws.Range("A436").Value = "DSS1435"
$ws.Range("B436").Value = "Mohamed Khairy Elsayed Elmasry"
$ws.Range("C436").Value = "Scaffold Competent Person"
$ws.Range("D436").Value = 45566
$ws.Range("E436").Value = 1

$ws.Range("A437").Value = "DSS1436"
$ws.Range("B437").Value = "Mohamed Khairy Elsayed Elmasry"
$ws.Range("C437").Value = "Lifting & Rigging Competent Person"
$ws.Range("D437").Value = 45566
$ws.Range("E437").Value = 1

$ws.Range("A438").Value = "DSS1437"
$ws.Range("B438").Value = "Mohamed Khairy Elsayed Elmasry"
$ws.Range("C438").Value = "Health & Safety Risk Assessment"
$ws.Range("D438").Value = 45566
$ws.Range("E438").Value = 1

$ws.Range("A439").Value = "DSS1438"
$ws.Range("B439").Value = "Mohamed Khairy Elsayed Elmasry"
$ws.Range("C439").Value = "Safety Management System & PTW"
$ws.Range("D439").Value = 45566
$ws.Range("E439").Value = 1

# ---------------------------------------------------------------------
# 3) Move the active selection the way the author left it.
# ---------------------------------------------------------------------
$ws.Range("A439").Select()
